$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first data row ("H 72") was removed; every subsequent row shifts up by one.
$ws.Rows.Item(2).Delete()
